# "Actualizacion de datadriven para estabilizar"
# Updates the ProyectoSolucion / ProyectoSolucion300 / Express300 / Express /
# PORTUGUES sheets of the datadriven workbook: fixes localisation strings,
# re-balances the "capacity %" test rows (adds a 20% row), flags the touched
# K-column cells in red, and tidies up a couple of sheet selections/column
# widths left over from manual testing.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "ProyectoSolucion"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("ProyectoSolucion")

# Row 2: only the K2 cell gets flagged red (value/content unchanged)
$ws1.Cells.Item(2, 11).Interior.Color = 255

# Row 3: several cells get new content, K3 also flagged red
$ws1.Range("D3").Value = "PRUEBA Auto70.1 Solucion Ganar Oferta"
$ws1.Range("G3").Value = "60%"
$ws1.Range("K3").Value = 70
$ws1.Cells.Item(3, 11).Interior.Color = 255
$ws1.Range("U3").ClearContents()
$ws1.Range("W3").Value = "100000003-Swaps"
$ws1.Range("Z3").Value = "Ganar"
$ws1.Range("AC3").Value = "C:/Users/Administrador/Documents/COLOMBIA-ProyectoINTERNEXA/Documento de orden de servicio Prueba Auto.txt"

# New row 4 (duplicate of the old row 3 pattern, with a 20% capacity case)
$ws1.Range("A4").Value = "jbedoya@internexa.com"
$ws1.Range("B4").Value = "Abr2018*"
$ws1.Range("C4").Value = "Avantel S.A."
$ws1.Range("D4").Value = "PRUEBA Auto37.1 Proyecto Perder Oferta"
$ws1.Range("E4").Value = "Nuevo cliente"
$ws1.Range("F4").Value = "23/04/2018"
$ws1.Range("G4").Value = "20%"
$ws1.Range("H4").Value = "24/04/2018"
$ws1.Range("I4").Value = "Proyecto"
$ws1.Range("J4").Value = 9
$ws1.Range("K4").Value = 37
$ws1.Cells.Item(4, 11).Interior.Color = 255
$ws1.Range("L4").Value = 120
$ws1.Range("M4").Value = "LIMA"
$ws1.Range("N4").Value = "bogotá"
$ws1.Range("O4").Value = "Av Republica de Panama 1123"
$ws1.Range("P4").Value = "Gelly Andrea Bustamante"
$ws1.Range("Q4").Value = "Pruebas factibilidad 12233"
$ws1.Range("R4").Value = "GBUSTAMANTE@INTERNEXA.COM"
$ws1.Range("S4").Value = "Abr2018*"
$ws1.Range("T4").Value = "ICA"
$ws1.Range("U4").Value = "CHICLAYO"
$ws1.Range("V4").Value = "FERNANDO PARRA"
$ws1.Range("W4").Value = "100000000-Cargo fijo mensual"
$ws1.Range("X4").Value = "100000000-Anticipado"
$ws1.Range("Y4").Value = "Nuqui"
$ws1.Range("Z4").Value = "Perder"
$ws1.Range("AA4").Value = 4
$ws1.Range("AB4").Value = "Ofera 12312"

# Column E got a bit wider (bestFit turned off in the process)
$ws1.Columns.Item(5).ColumnWidth = 13.5

$ws1.Activate()
$ws1.Range("E2:E4").Select()

# ---------------------------------------------------------------------
# Sheet "ProyectoSolucion300"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ProyectoSolucion300")
$ws2.Activate()
$ws2.Range("E2").Select()

# ---------------------------------------------------------------------
# Sheet "Express300"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Express300")
$ws3.Range("E2").Value = "Novo cliente"
$ws3.Range("E3").Value = "Novo cliente"
$ws3.Activate()
$ws3.Range("D23").Select()

# ---------------------------------------------------------------------
# Sheet "Express"
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Express")
$ws4.Range("E2").Value = "Novo cliente"
$ws4.Range("E3").Value = "Novo cliente"
$ws4.Range("E4").Value = "Novo cliente"
$ws4.Range("E5").Value = "Novo cliente"
$ws4.Range("E6").Value = "Novo cliente"
$ws4.Range("E7").Value = "Novo cliente"
$ws4.Activate()
$ws4.Range("E10").Select()

# ---------------------------------------------------------------------
# Sheet "PORTUGUES"
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("PORTUGUES")
$ws6.Columns.Item(4).ColumnWidth = 29.6
$ws6.Columns.Item(5).ColumnWidth = 14.6
$ws6.Activate()
$ws6.Range("A19").Select()
$ws6.Application.ActiveWindow.ScrollRow = 19
$ws6.Range("A24:XFD25,A28:XFD28").Select()
$ws6.Range("A28").Activate()

# ---------------------------------------------------------------------
# Restore "ProyectoSolucion" as the active/visible tab
# ---------------------------------------------------------------------
$ws1.Activate()
